# Weekly cryptos data refresh: update Price (D) and Volume(1h) (E) columns,
# and fix the swapped Kaspa / WEMIXToken rows (37-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.222.04'
$ws.Range("E2").Value = '  +0.88%  '

# Row 3
$ws.Range("D3").Value = '2.271.99'
$ws.Range("E3").Value = '  +0.01%  '

# Row 4
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.93'
$ws.Range("E5").Value = '  +0.85%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.64'
$ws.Range("E6").Value = '  +5.18%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +1.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.64'
$ws.Range("E10").Value = '  +9.47%  '

# Row 12
$ws.Range("E12").Value = '  -0.93%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.67'
$ws.Range("E13").Value = '  -0.44%  '

# Row 14
$ws.Range("D14").Value = '2.595.71'
$ws.Range("E14").Value = '  -1.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.44'
$ws.Range("E15").Value = '  +0.89%  '

# Row 16
$ws.Range("D16").Value = '2.266.69'
$ws.Range("E16").Value = '  -0.40%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.797'
$ws.Range("E17").Value = '  +2.05%  '

# Row 18
$ws.Range("D18").Value = '42.126.22'
$ws.Range("E18").Value = '  +0.85%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.51'
$ws.Range("E19").Value = '  -2.26%  '

# Row 20
$ws.Range("E20").Value = '  +0.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.99'
$ws.Range("E21").Value = '  +0.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.80'
$ws.Range("E22").Value = '  +0.71%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.23'
$ws.Range("E23").Value = '  -2.46%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  +0.29%  '

# Row 25
$ws.Range("E25").Value = '  +0.84%  '

# Row 26
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.76'
$ws.Range("E27").Value = '  -1.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.40'
$ws.Range("E28").Value = '  +6.93%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.53'
$ws.Range("E29").Value = '  -0.25%  '

# Row 30
$ws.Range("E30").Value = '  +1.66%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.39'
$ws.Range("E31").Value = '  -0.06%  '

# Row 32
$ws.Range("E32").Value = '  +0.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.17'
$ws.Range("E34").Value = '  +4.65%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0741'
$ws.Range("E35").Value = '  -0.28%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.30'
$ws.Range("E36").Value = '  +2.53%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.105'
$ws.Range("E37").Value = '  -0.46%  '

# Row 38
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -0.31%  '

# Row 39
$ws.Range("E39").Value = '  +2.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.115'
$ws.Range("E40").Value = '  -1.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.09'
$ws.Range("E41").Value = '  +3.56%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.45'
$ws.Range("E42").Value = '  +14.46%  '

# Row 43
$ws.Range("D43").Value = '1.987.55'
$ws.Range("E43").Value = '  -1.33%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0284'
$ws.Range("E44").Value = '  +0.83%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.98'
$ws.Range("E45").Value = '  -2.40%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.94'
$ws.Range("E46").Value = '  +1.37%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.94'
$ws.Range("E47").Value = '  -4.85%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.26'
$ws.Range("E48").Value = '  -0.10%  '

# Row 49
$ws.Range("E49").Value = '  +1.30%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.05'
$ws.Range("E50").Value = '  -1.51%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '91.64'
$ws.Range("E51").Value = '  -0.22%  '
